# CHANGED: Card List - Increased threat of player effects and power of non-unique enemy starships.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Increase Threat (column E) of player effects by 1
$threatRows = @(10, 11, 12, 13, 42, 43, 45, 59, 60, 62)
foreach ($r in $threatRows) {
    $cell = $ws.Cells.Item($r, 5)   # column E = Threat
    $cell.Value2 = $cell.Value2 + 1
}

# Increase Power (column F) of non-unique enemy starships by 1
$powerRows = @(26, 27, 35)
foreach ($r in $powerRows) {
    $cell = $ws.Cells.Item($r, 6)   # column F = Power
    $cell.Value2 = $cell.Value2 + 1
}

# Add extra game text for the Swarm Frigate (row 35)
$ws.Range("K35").Value = "Add (2)."

# Update the saved view state (scroll position / active selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B25").Select()
